# mentoring.xlsx -- "Total revamp of project" update
#
# Adds newly-finished/new mentees across sheets (Graduate Students,
# Postdoctoral Fellows, Junior Faculty), fills in a couple of "current
# position" cells for people who have since moved on, and updates the
# row-selection / active-sheet UI state to match.

$wb = $excel.ActiveWorkbook

$wsGrad    = $wb.Worksheets.Item("Graduate Students")
$wsPostdoc = $wb.Worksheets.Item("Postdoctoral Fellows")
$wsFaculty = $wb.Worksheets.Item("Junior Faculty")
$wsOther   = $wb.Worksheets.Item("Other Trainees")

# ---------------------------------------------------------------
# Graduate Students: add Ziyou Ren as row 25
# ---------------------------------------------------------------
$wsGrad.Range("A24").Copy($wsGrad.Range("A25"))
$wsGrad.Range("A25").Formula = "=A24+1"
$wsGrad.Range("B25").Value = 2017
$wsGrad.Range("C25").Value = "present"
$wsGrad.Range("D25").Value = "Ziyou Ren"

# ---------------------------------------------------------------
# Junior Faculty: Curtis Weiss is now ICU Director at Northshore Medicine
# ---------------------------------------------------------------
$wsFaculty.Range("H4").Value = "Northshore Medicine"
$wsFaculty.Range("G4").Value = "ICU Director"

# ---------------------------------------------------------------
# back to Graduate Students: finish Ziyou Ren's row
# ---------------------------------------------------------------
$wsGrad.Range("E25").Value = "co-advised with Scott Buddinger"
$wsGrad.Range("F25").Value = "-"
$wsGrad.Range("G25").Value = "-"
$wsGrad.Range("H25").Value = "-"

# ---------------------------------------------------------------
# Postdoctoral Fellows: Diego is now a Research Scientist at the
# Army Research Laboratory
# ---------------------------------------------------------------
$wsPostdoc.Range("G16").Value = "Research Scientist"
$wsPostdoc.Range("H16").Value = "Army Research Laboratory"

# ---------------------------------------------------------------
# Junior Faculty: add Paul A. Reyfman, M.D. as row 6 (and convert
# Adam Pah's row-5 "A" counter into the same kind of running formula)
# ---------------------------------------------------------------
$wsFaculty.Range("A5").Formula = "=A4+1"
$wsFaculty.Range("A5:H5").Copy($wsFaculty.Range("A6:H6"))
$wsFaculty.Range("A6").Formula = "=A5+1"
$wsFaculty.Range("B6").Value = 2017
$wsFaculty.Range("C6").Value = "present"
$wsFaculty.Range("D6").Value = "Paul A. Reyfman, M.D."
$wsFaculty.Range("E6").Value = "Instructor -- F32 award"
$wsFaculty.Range("F6").Value = "-"
$wsFaculty.Range("G6").Value = "-"
$wsFaculty.Range("H6").Value = "Northwestern University"

# ---------------------------------------------------------------
# Postdoctoral Fellows: Yang Yang's postdoc ended in 2017 and she is
# now a Data Scientist
# ---------------------------------------------------------------
$wsPostdoc.Range("C19").Value = 2017
$wsPostdoc.Range("G19").Value = "Data Scientist"

# ---------------------------------------------------------------
# Postdoctoral Fellows: add Julia Poncela-Casasnovas as row 20, plus
# three trailing (still empty) styled rows
# ---------------------------------------------------------------
$wsPostdoc.Range("A18").Copy($wsPostdoc.Range("A20"))
$wsPostdoc.Range("A20").Formula = "=A19+1"
$wsPostdoc.Range("B20").Value = 2018
$wsPostdoc.Range("C20").Value = "present"
$wsPostdoc.Range("D20").Value = "Julia Poncela-Casasnovas"
$wsPostdoc.Range("E20").Value = "-"
$wsPostdoc.Range("F19").Copy($wsPostdoc.Range("F20"))
$wsPostdoc.Range("F20").Value = "-"
$wsPostdoc.Range("G20").Value = "Postdoctoral Fellow"
$wsPostdoc.Range("H20").Value = "Northwestern University"

$wsPostdoc.Range("A19").Copy($wsPostdoc.Range("A21"))
$wsPostdoc.Range("A21").ClearContents()
$wsPostdoc.Range("A19").Copy($wsPostdoc.Range("A22"))
$wsPostdoc.Range("A22").ClearContents()
$wsPostdoc.Range("A19").Copy($wsPostdoc.Range("A23"))
$wsPostdoc.Range("A23").ClearContents()

# ---------------------------------------------------------------
# Selections: restore each sheet's last-used cell; leave the
# "Postdoctoral Fellows" tab active (it was "Other Trainees" before)
# ---------------------------------------------------------------
$wsGrad.Range("F25").Select()
$wsFaculty.Range("E7").Select()
$wsOther.Range("F55").Select()
$wsPostdoc.Range("F20").Select()

Write-Output "mentoring.xlsx update applied"
